$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 72836
$ws.Range("E2").Value = 1496
$ws.Range("F2").Value = 1496
$ws.Range("G2").Value = 342
$ws.Range("H2").Value = 140
$ws.Range("I2").Value = -251
$ws.Range("J2").Value = 391
$ws.Range("K2").Value = 55548
$ws.Range("L2").Value = 35487
$ws.Range("M2").Value = 20061
$ws.Range("N2").Value = 12278
$ws.Range("O2").Value = 7783
$ws.Range("P2").Value = 1183
$ws.Range("Q2").Value = 1864
$ws.Range("R2").Value = -8976
$ws.Range("S2").Value = 5259
$ws.Range("T2").Value = 6595
$ws.Range("U2").Value = -4731
$ws.Range("V2").Value = 27552
$ws.Range("W2").Value = 2.05
$ws.Range("X2").Value = 0.19
$ws.Range("Y2").Value = -2.01
$ws.Range("Z2").Value = 0.27
$ws.Range("AA2").Value = 176.89
$ws.Range("AB2").Value = 998.64
$ws.Range("AC2").Value = -1032
$ws.Range("AD2").Value = -74.65000000000001
$ws.Range("AE2").Value = 58227
$ws.Range("AF2").Value = 1.32
$ws.Range("AG2").Value = 291
$ws.Range("AH2").Value = 0.38
$ws.Range("AI2").Value = -24.98
$ws.Range("AJ2").Value = 21505409

# --- Row 3 ---
$ws.Range("D3").Value = 52692
$ws.Range("E3").Value = 1199
$ws.Range("F3").Value = 1143
$ws.Range("G3").Value = 1195
$ws.Range("H3").Value = 944
$ws.Range("I3").Value = 556
$ws.Range("J3").Value = 388
$ws.Range("K3").Value = 65950
$ws.Range("L3").Value = 40092
$ws.Range("M3").Value = 25858
$ws.Range("N3").Value = 14975
$ws.Range("O3").Value = 10882
$ws.Range("P3").Value = 1356
$ws.Range("Q3").Value = -654
$ws.Range("R3").Value = -6594
$ws.Range("S3").Value = 9560
$ws.Range("T3").Value = 6671
$ws.Range("U3").Value = -7325
$ws.Range("V3").Value = 32406
$ws.Range("W3").Value = 2.28
$ws.Range("X3").Value = 1.79
$ws.Range("Y3").Value = 4.08
$ws.Range("Z3").Value = 1.55
$ws.Range("AA3").Value = 155.05
$ws.Range("AB3").Value = 1053.59
$ws.Range("AC3").Value = 2269
$ws.Range("AD3").Value = 39.51
$ws.Range("AE3").Value = 62701
$ws.Range("AF3").Value = 1.43
$ws.Range("AG3").Value = 300
$ws.Range("AH3").Value = 0.33
$ws.Range("AI3").Value = 13.15
$ws.Range("AJ3").Value = 24303255

# --- Row 4 ---
$ws.Range("D4").Value = 53222
$ws.Range("E4").Value = 1897
$ws.Range("F4").Value = 2304
$ws.Range("G4").Value = 1958
$ws.Range("H4").Value = 1759
$ws.Range("I4").Value = 826
$ws.Range("J4").Value = 932
$ws.Range("K4").Value = 65459
$ws.Range("L4").Value = 38835
$ws.Range("M4").Value = 26624
$ws.Range("N4").Value = 15893
$ws.Range("O4").Value = 10732
$ws.Range("P4").Value = 1356
$ws.Range("Q4").Value = -325
$ws.Range("R4").Value = -3259
$ws.Range("S4").Value = 1725
$ws.Range("T4").Value = 3473
$ws.Range("U4").Value = -3798
$ws.Range("V4").Value = 27513
$ws.Range("W4").Value = 3.56
$ws.Range("X4").Value = 3.3
$ws.Range("Y4").Value = 5.35
$ws.Range("Z4").Value = 2.68
$ws.Range("AA4").Value = 145.86
$ws.Range("AB4").Value = 1106.27
$ws.Range("AC4").Value = 3046
$ws.Range("AD4").Value = 27.17
$ws.Range("AE4").Value = 66542
$ws.Range("AF4").Value = 1.24
$ws.Range("AG4").Value = 350
$ws.Range("AH4").Value = 0.42
$ws.Range("AI4").Value = 10.29
$ws.Range("AJ4").Value = 24303255

# --- Row 5 ---
$ws.Range("D5").Value = 64159
$ws.Range("E5").Value = 821
$ws.Range("F5").Value = 821
$ws.Range("G5").Value = 1511
$ws.Range("H5").Value = 8012
$ws.Range("I5").Value = 7179
$ws.Range("J5").Value = 833
$ws.Range("K5").Value = 50516
$ws.Range("L5").Value = 29080
$ws.Range("M5").Value = 21436
$ws.Range("N5").Value = 10936
$ws.Range("O5").Value = 10500
$ws.Range("P5").Value = 704
$ws.Range("Q5").Value = -1306
$ws.Range("R5").Value = 72
$ws.Range("S5").Value = 513
$ws.Range("T5").Value = 1941
$ws.Range("U5").Value = -3247
$ws.Range("V5").Value = 21350
$ws.Range("W5").Value = 1.28
$ws.Range("X5").Value = 12.49
$ws.Range("Y5").Value = 53.52
$ws.Range("Z5").Value = 13.82
$ws.Range("AA5").Value = 135.66
$ws.Range("AB5").Value = 2934.94
$ws.Range("AC5").Value = 28731
$ws.Range("AD5").Value = 3.62
$ws.Range("AE5").Value = 90161
$ws.Range("AF5").Value = 1.15
$ws.Range("AG5").Value = 750
$ws.Range("AH5").Value = 0.72
$ws.Range("AI5").Value = 1.28
$ws.Range("AJ5").Value = 10784283

# --- Row 6 ---
$ws.Range("D6").Value = 69392
$ws.Range("E6").Value = 1000
$ws.Range("F6").Value = 1000
$ws.Range("G6").Value = 701
$ws.Range("H6").Value = 977
$ws.Range("I6").Value = 708
$ws.Range("K6").Value = 48056
$ws.Range("L6").Value = 24870
$ws.Range("M6").Value = 23186
$ws.Range("N6").Value = 14495
$ws.Range("P6").Value = 1117
$ws.Range("Q6").Value = -1581
$ws.Range("R6").Value = -1594
$ws.Range("S6").Value = 2993
$ws.Range("T6").Value = 2558
$ws.Range("U6").Value = -4140
$ws.Range("V6").Value = 15211
$ws.Range("W6").Value = 1.44
$ws.Range("X6").Value = 1.41
$ws.Range("Y6").Value = 5.56
$ws.Range("Z6").Value = 1.98
$ws.Range("AA6").Value = 107.27
$ws.Range("AB6").Value = 2118.77
$ws.Range("AC6").Value = 3936
$ws.Range("AD6").Value = 6.63
$ws.Range("AE6").Value = 71115
$ws.Range("AF6").Value = 0.37
$ws.Range("AI6").Value = 17.38
$ws.Range("AJ6").Value = 19037855

# AG6 and AH6 no longer exist after the correction -> clear them entirely
$ws.Range("AG6:AH6").ClearContents()

# Rows 7-9: these data rows only keep columns A-C (index/ticker/company); all financial
# figures (D:AJ) are removed as part of the error fix
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
